$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 7423831.5
$ws.Cells.Item(62, 9).Value = 13897309
$ws.Cells.Item(62, 10).Value = 25571.428
$ws.Cells.Item(62, 11).Value = 13897309
$ws.Cells.Item(62, 12).Value = 25571.428
$ws.Cells.Item(62, 13).Value = -13896685
$ws.Cells.Item(62, 14).Value = -26819.428

$ws.Cells.Item(65, 8).Value = 7423831.5
$ws.Cells.Item(65, 9).Value = 13897309
$ws.Cells.Item(65, 10).Value = 25571.428
$ws.Cells.Item(65, 11).Value = 69486545
$ws.Cells.Item(65, 12).Value = 127857.14
$ws.Cells.Item(65, 13).Value = -69483425
$ws.Cells.Item(65, 14).Value = -134097.14

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 127113.875
$ws.Cells.Item(2, 9).Value = 145130.14
$ws.Cells.Item(2, 10).Value = 1000
$ws.Cells.Item(2, 11).Value = 145130.14
$ws.Cells.Item(2, 12).Value = 1000
$ws.Cells.Item(2, 13).Value = -145017.14
$ws.Cells.Item(2, 14).Value = -1226

$ws.Cells.Item(32, 8).Value = 17133.139
$ws.Cells.Item(32, 9).Value = 4432.7593
$ws.Cells.Item(32, 11).Value = 4432.7593
$ws.Cells.Item(32, 13).Value = -4145.7593

$ws.Cells.Item(61, 8).Value = 2700.1365
$ws.Cells.Item(61, 9).Value = 1855.0625
$ws.Cells.Item(61, 10).Value = 4953.6665
$ws.Cells.Item(61, 11).Value = 1855.0625
$ws.Cells.Item(61, 12).Value = 4953.6665
$ws.Cells.Item(61, 13).Value = -1643.0625
$ws.Cells.Item(61, 14).Value = -5377.6665

$ws.Cells.Item(102, 8).Value = 3011.7334
$ws.Cells.Item(102, 9).Value = 3370.625
$ws.Cells.Item(102, 10).Value = 2601.5715
$ws.Cells.Item(102, 11).Value = 3370.625
$ws.Cells.Item(102, 12).Value = 2601.5715
$ws.Cells.Item(102, 13).Value = -1748.625
$ws.Cells.Item(102, 14).Value = -5845.5715

$ws.Cells.Item(116, 8).Value = 127113.875
$ws.Cells.Item(116, 9).Value = 145130.14
$ws.Cells.Item(116, 10).Value = 1000
$ws.Cells.Item(116, 11).Value = 145130.14
$ws.Cells.Item(116, 12).Value = 1000
$ws.Cells.Item(116, 13).Value = -142836.14
$ws.Cells.Item(116, 14).Value = -5588

$ws.Cells.Item(122, 8).Value = 1839.5555
$ws.Cells.Item(122, 9).Value = 1573.6
$ws.Cells.Item(122, 10).Value = 2172
$ws.Cells.Item(122, 11).Value = 4720.799999999999
$ws.Cells.Item(122, 12).Value = 6516
$ws.Cells.Item(122, 13).Value = -2270.799999999999
$ws.Cells.Item(122, 14).Value = -11416

$ws.Cells.Item(132, 8).Value = 2621.5
$ws.Cells.Item(132, 9).Value = 2070.9534
$ws.Cells.Item(132, 10).Value = 4773.636
$ws.Cells.Item(132, 11).Value = 6212.860199999999
$ws.Cells.Item(132, 12).Value = 14320.908
$ws.Cells.Item(132, 13).Value = -3682.860199999999
$ws.Cells.Item(132, 14).Value = -19380.908

$ws.Cells.Item(136, 8).Value = 2700.1365
$ws.Cells.Item(136, 9).Value = 1855.0625
$ws.Cells.Item(136, 10).Value = 4953.6665
$ws.Cells.Item(136, 11).Value = 5565.1875
$ws.Cells.Item(136, 12).Value = 14860.9995
$ws.Cells.Item(136, 13).Value = -3015.1875
$ws.Cells.Item(136, 14).Value = -19960.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 127113.875
$ws.Cells.Item(3, 9).Value = 145130.14
$ws.Cells.Item(3, 10).Value = 1000
$ws.Cells.Item(3, 11).Value = 145130.14
$ws.Cells.Item(3, 12).Value = 1000
$ws.Cells.Item(3, 13).Value = -145016.14
$ws.Cells.Item(3, 14).Value = -1228

$ws.Cells.Item(134, 8).Value = 2825.6365
$ws.Cells.Item(134, 9).Value = 1933.8857
$ws.Cells.Item(134, 11).Value = 5801.6571
$ws.Cells.Item(134, 13).Value = -3266.6571

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 891.73334
$ws.Cells.Item(22, 9).Value = 592
$ws.Cells.Item(22, 10).Value = 1341.3334
$ws.Cells.Item(22, 11).Value = 592
$ws.Cells.Item(22, 12).Value = 1341.3334
$ws.Cells.Item(22, 13).Value = -242
$ws.Cells.Item(22, 14).Value = -2041.3334

$ws.Cells.Item(58, 8).Value = 37039010
$ws.Cells.Item(58, 9).Value = 47620290
$ws.Cells.Item(58, 10).Value = 4518.5
$ws.Cells.Item(58, 11).Value = 47620290
$ws.Cells.Item(58, 12).Value = 4518.5
$ws.Cells.Item(58, 13).Value = -47620087
$ws.Cells.Item(58, 14).Value = -4924.5

$ws.Cells.Item(107, 8).Value = 447.77777
$ws.Cells.Item(107, 9).Value = 336.41666
$ws.Cells.Item(107, 11).Value = 336.41666
$ws.Cells.Item(107, 13).Value = 1583.58334

$ws.Cells.Item(132, 8).Value = 10419620
$ws.Cells.Item(132, 9).Value = 20835096
$ws.Cells.Item(132, 11).Value = 62505288
$ws.Cells.Item(132, 13).Value = -62502758

$ws.Cells.Item(134, 8).Value = 30614470
$ws.Cells.Item(134, 9).Value = 31251280
$ws.Cells.Item(134, 10).Value = 29415772
$ws.Cells.Item(134, 11).Value = 93753840
$ws.Cells.Item(134, 12).Value = 88247316
$ws.Cells.Item(134, 13).Value = -93751305
$ws.Cells.Item(134, 14).Value = -88252386

$ws.Cells.Item(136, 8).Value = 37039010
$ws.Cells.Item(136, 9).Value = 47620290
$ws.Cells.Item(136, 10).Value = 4518.5
$ws.Cells.Item(136, 11).Value = 142860870
$ws.Cells.Item(136, 12).Value = 13555.5
$ws.Cells.Item(136, 13).Value = -142858320
$ws.Cells.Item(136, 14).Value = -18655.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 4628.2144
$ws.Cells.Item(34, 10).Value = 4953.4614
$ws.Cells.Item(34, 12).Value = 14860.3842
$ws.Cells.Item(34, 14).Value = -15028.3842

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 5250.2563
$ws.Cells.Item(70, 9).Value = 5350.645
$ws.Cells.Item(70, 10).Value = 4861.25
$ws.Cells.Item(70, 11).Value = 5350.645
$ws.Cells.Item(70, 12).Value = 4861.25
$ws.Cells.Item(70, 13).Value = -5080.645
$ws.Cells.Item(70, 14).Value = -5401.25

$ws.Cells.Item(73, 8).Value = 5250.2563
$ws.Cells.Item(73, 9).Value = 5350.645
$ws.Cells.Item(73, 10).Value = 4861.25
$ws.Cells.Item(73, 11).Value = 5350.645
$ws.Cells.Item(73, 12).Value = 4861.25
$ws.Cells.Item(73, 13).Value = -4414.645
$ws.Cells.Item(73, 14).Value = -6733.25

$ws.Cells.Item(122, 8).Value = 1837.875
$ws.Cells.Item(122, 9).Value = 1600.4286
$ws.Cells.Item(122, 10).Value = 3500
$ws.Cells.Item(122, 11).Value = 4801.2858
$ws.Cells.Item(122, 12).Value = 10500
$ws.Cells.Item(122, 13).Value = -2351.2858
$ws.Cells.Item(122, 14).Value = -15400

$ws.Cells.Item(132, 8).Value = 3550.8518
$ws.Cells.Item(132, 9).Value = 3478.96
$ws.Cells.Item(132, 11).Value = 10436.88
$ws.Cells.Item(132, 13).Value = -7906.880000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(32, 8).Value = 1156.5
$ws.Cells.Item(32, 9).Value = 1156.5
$ws.Cells.Item(32, 11).Value = 1156.5
$ws.Cells.Item(32, 13).Value = -839.5

$ws.Cells.Item(122, 8).Value = 3531.6785
$ws.Cells.Item(122, 9).Value = 2954.2727
$ws.Cells.Item(122, 10).Value = 3905.2942
$ws.Cells.Item(122, 11).Value = 8862.8181
$ws.Cells.Item(122, 12).Value = 11715.8826
$ws.Cells.Item(122, 13).Value = -6412.8181
$ws.Cells.Item(122, 14).Value = -16615.8826

$ws.Cells.Item(132, 8).Value = 4557.85
$ws.Cells.Item(132, 9).Value = 2651.4
$ws.Cells.Item(132, 10).Value = 6464.3
$ws.Cells.Item(132, 11).Value = 7954.200000000001
$ws.Cells.Item(132, 12).Value = 19392.9
$ws.Cells.Item(132, 13).Value = -5424.200000000001
$ws.Cells.Item(132, 14).Value = -24452.9

$ws.Cells.Item(136, 8).Value = 5165.933
$ws.Cells.Item(136, 9).Value = 2629.4092
$ws.Cells.Item(136, 11).Value = 7888.2276
$ws.Cells.Item(136, 13).Value = -5338.2276

$ws.Cells.Item(139, 8).Value = 53850
$ws.Cells.Item(139, 10).Value = 53850
$ws.Cells.Item(139, 12).Value = 53850
$ws.Cells.Item(139, 14).Value = -64130

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 2290.8508
$ws.Cells.Item(132, 9).Value = 2125.6123
$ws.Cells.Item(132, 11).Value = 6376.836899999999
$ws.Cells.Item(132, 13).Value = -3846.836899999999

$ws.Cells.Item(136, 8).Value = 2912.5715
$ws.Cells.Item(136, 9).Value = 1091.1538
$ws.Cells.Item(136, 11).Value = 3273.4614
$ws.Cells.Item(136, 13).Value = -723.4614000000001
